$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 11 - BLEU score
$ws.Range("B11").Value = 0.1267333857233461

# Row 12 - Code BLEU
$ws.Range("B12").Value = 0.3257498546112294
$ws.Range("C12").Value = "{'codebleu': 0.3257498546112294, 'ngram_match_score': 0.12640884173887587, 'weighted_ngram_match_score': 0.14447270517752783, 'syntax_match_score': 0.5377358490566038, 'dataflow_match_score': 0.4943820224719101}"

# Row 13 - Embeddings and Cosine similarity
$ws.Range("B13").Value = 0.8913486560164019
